$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 92
Set-CellValue $ws 92 8 741517.9399999999  # H92: 556061.9 -> 741517.9399999999
Set-CellValue $ws 92 9 741517.9399999999  # I92: 585307.2 -> 741517.9399999999
Set-CellValue $ws 92 10 0  # J92: 400 -> 0
Set-CellValue $ws 92 11 741517.9399999999  # K92: 585307.2 -> 741517.9399999999
Set-CellValue $ws 92 12 0  # L92: 400 -> 0
Set-CellValue $ws 92 13 -740269.9399999999  # M92: -584059.2 -> -740269.9399999999
Set-CellValue $ws 92 14 $null  # N92: -2896 -> (removed)
# Row 133
Set-CellValue $ws 133 8 10797  # H133: 14791.9 -> 10797
Set-CellValue $ws 133 10 10797  # J133: 14791.9 -> 10797
Set-CellValue $ws 133 12 10797  # L133: 14791.9 -> 10797
Set-CellValue $ws 133 14 -20917  # N133: -24911.9 -> -20917
# Row 135
Set-CellValue $ws 135 8 1746.88  # H135: 1517.3115 -> 1746.88
Set-CellValue $ws 135 9 1680.2858  # I135: 1422.5927 -> 1680.2858
Set-CellValue $ws 135 10 2096.5  # J135: 2248 -> 2096.5
Set-CellValue $ws 135 11 15122.5722  # K135: 12803.3343 -> 15122.5722
Set-CellValue $ws 135 12 18868.5  # L135: 20232 -> 18868.5
Set-CellValue $ws 135 13 -12587.5722  # M135: -10268.3343 -> -12587.5722
Set-CellValue $ws 135 14 -23938.5  # N135: -25302 -> -23938.5
# Row 138
Set-CellValue $ws 138 8 4119468.8  # H138: 4119578 -> 4119468.8
Set-CellValue $ws 138 9 1167614.5  # I138: 1468844.5 -> 1167614.5
Set-CellValue $ws 138 10 6291588  # J138: 5466672 -> 6291588
Set-CellValue $ws 138 11 3502843.5  # K138: 4406533.5 -> 3502843.5
Set-CellValue $ws 138 12 18874764  # L138: 16400016 -> 18874764
Set-CellValue $ws 138 13 -3497703.5  # M138: -4401393.5 -> -3497703.5
Set-CellValue $ws 138 14 -18885044  # N138: -16410296 -> -18885044
# Row 141
Set-CellValue $ws 141 8 2080.9722  # H141: 2652.2373 -> 2080.9722
Set-CellValue $ws 141 9 1399.1428  # I141: 1724.9149 -> 1399.1428
Set-CellValue $ws 141 10 6853.778  # J141: 6284.25 -> 6853.778
Set-CellValue $ws 141 11 4197.428400000001  # K141: 5174.7447 -> 4197.428400000001
Set-CellValue $ws 141 12 20561.334  # L141: 18852.75 -> 20561.334
Set-CellValue $ws 141 13 982.5715999999993  # M141: 5.255299999999806 -> 982.5715999999993
Set-CellValue $ws 141 14 -30921.334  # N141: -29212.75 -> -30921.334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
Set-CellValue $ws 32 8 13779.033  # H32: 15285.524 -> 13779.033
Set-CellValue $ws 32 9 2550.0488  # I32: 2813.7534 -> 2550.0488
Set-CellValue $ws 32 10 128876.125  # J32: 116445.445 -> 128876.125
Set-CellValue $ws 32 11 2550.0488  # K32: 2813.7534 -> 2550.0488
Set-CellValue $ws 32 12 128876.125  # L32: 116445.445 -> 128876.125
Set-CellValue $ws 32 13 -2263.0488  # M32: -2526.7534 -> -2263.0488
Set-CellValue $ws 32 14 -129450.125  # N32: -117019.445 -> -129450.125
# Row 61
Set-CellValue $ws 61 8 1977.4706  # H61: 1857.9636 -> 1977.4706
Set-CellValue $ws 61 9 1562.9048  # I61: 1501.4773 -> 1562.9048
Set-CellValue $ws 61 10 3912.111  # J61: 3283.9092 -> 3912.111
Set-CellValue $ws 61 11 1562.9048  # K61: 1501.4773 -> 1562.9048
Set-CellValue $ws 61 12 3912.111  # L61: 3283.9092 -> 3912.111
Set-CellValue $ws 61 13 -1350.9048  # M61: -1289.4773 -> -1350.9048
Set-CellValue $ws 61 14 -4336.111  # N61: -3707.9092 -> -4336.111
# Row 74
Set-CellValue $ws 74 8 9841.177  # H74: 6594.08 -> 9841.177
Set-CellValue $ws 74 9 2588.889  # I74: 1240.75 -> 2588.889
Set-CellValue $ws 74 10 18000  # J74: 16111.111 -> 18000
Set-CellValue $ws 74 11 2588.889  # K74: 1240.75 -> 2588.889
Set-CellValue $ws 74 12 18000  # L74: 16111.111 -> 18000
Set-CellValue $ws 74 13 -1714.889  # M74: -366.75 -> -1714.889
Set-CellValue $ws 74 14 -19748  # N74: -17859.111 -> -19748
# Row 77
Set-CellValue $ws 77 8 9841.177  # H77: 6594.08 -> 9841.177
Set-CellValue $ws 77 9 2588.889  # I77: 1240.75 -> 2588.889
Set-CellValue $ws 77 10 18000  # J77: 16111.111 -> 18000
Set-CellValue $ws 77 11 12944.445  # K77: 6203.75 -> 12944.445
Set-CellValue $ws 77 12 90000  # L77: 80555.55500000001 -> 90000
Set-CellValue $ws 77 13 -8576.445  # M77: -1835.75 -> -8576.445
Set-CellValue $ws 77 14 -98736  # N77: -89291.55500000001 -> -98736
# Row 102
Set-CellValue $ws 102 8 3177  # H102: 3247.8948 -> 3177
Set-CellValue $ws 102 9 3819.182  # I102: 3593.1538 -> 3819.182
Set-CellValue $ws 102 10 1999.6666  # J102: 2499.8333 -> 1999.6666
Set-CellValue $ws 102 11 3819.182  # K102: 3593.1538 -> 3819.182
Set-CellValue $ws 102 12 1999.6666  # L102: 2499.8333 -> 1999.6666
Set-CellValue $ws 102 13 -2197.182  # M102: -1971.1538 -> -2197.182
Set-CellValue $ws 102 14 -5243.6666  # N102: -5743.8333 -> -5243.6666
# Row 122
Set-CellValue $ws 122 8 16360.429  # H122: 22100 -> 16360.429
Set-CellValue $ws 122 9 18670.5  # I122: 27000 -> 18670.5
Set-CellValue $ws 122 11 56011.5  # K122: 81000 -> 56011.5
Set-CellValue $ws 122 13 -53561.5  # M122: -78550 -> -53561.5
# Row 133
Set-CellValue $ws 133 8 52499.75  # H133: 51999.8 -> 52499.75
Set-CellValue $ws 133 10 52499.75  # J133: 51999.8 -> 52499.75
Set-CellValue $ws 133 12 52499.75  # L133: 51999.8 -> 52499.75
Set-CellValue $ws 133 14 -57559.75  # N133: -57059.8 -> -57559.75
# Row 136
Set-CellValue $ws 136 8 1977.4706  # H136: 1857.9636 -> 1977.4706
Set-CellValue $ws 136 9 1562.9048  # I136: 1501.4773 -> 1562.9048
Set-CellValue $ws 136 10 3912.111  # J136: 3283.9092 -> 3912.111
Set-CellValue $ws 136 11 4688.7144  # K136: 4504.4319 -> 4688.7144
Set-CellValue $ws 136 12 11736.333  # L136: 9851.7276 -> 11736.333
Set-CellValue $ws 136 13 -2138.7144  # M136: -1954.4319 -> -2138.7144
Set-CellValue $ws 136 14 -16836.333  # N136: -14951.7276 -> -16836.333
# Row 139
Set-CellValue $ws 139 8 42148.637  # H139: 46490.625 -> 42148.637
Set-CellValue $ws 139 9 30500  # I139: 0 -> 30500
Set-CellValue $ws 139 10 44737.223  # J139: 46490.625 -> 44737.223
Set-CellValue $ws 139 11 30500  # K139: 0 -> 30500
Set-CellValue $ws 139 12 44737.223  # L139: 46490.625 -> 44737.223
Set-CellValue $ws 139 13 -25360  # M139: None -> -25360
Set-CellValue $ws 139 14 -55017.223  # N139: -56770.625 -> -55017.223

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 59
Set-CellValue $ws 59 8 49999.5  # H59: 50000 -> 49999.5
Set-CellValue $ws 59 10 49999.5  # J59: 50000 -> 49999.5
Set-CellValue $ws 59 12 49999.5  # L59: 50000 -> 49999.5
Set-CellValue $ws 59 14 -51693.5  # N59: -51694 -> -51693.5
# Row 64
Set-CellValue $ws 64 8 618  # H64: 647.5 -> 618
Set-CellValue $ws 64 10 588.8570999999999  # J64: 624.4 -> 588.8570999999999
Set-CellValue $ws 64 12 588.8570999999999  # L64: 624.4 -> 588.8570999999999
Set-CellValue $ws 64 14 -1038.8571  # N64: -1074.4 -> -1038.8571
# Row 67
Set-CellValue $ws 67 8 618  # H67: 647.5 -> 618
Set-CellValue $ws 67 10 588.8570999999999  # J67: 624.4 -> 588.8570999999999
Set-CellValue $ws 67 12 588.8570999999999  # L67: 624.4 -> 588.8570999999999
Set-CellValue $ws 67 14 -2148.8571  # N67: -2184.4 -> -2148.8571
# Row 134
Set-CellValue $ws 134 8 27780744  # H134: 25643976 -> 27780744
Set-CellValue $ws 134 9 35716750  # I134: 35716880 -> 35716750
Set-CellValue $ws 134 10 4715.75  # J134: 3858.5454 -> 4715.75
Set-CellValue $ws 134 11 107150250  # K134: 107150640 -> 107150250
Set-CellValue $ws 134 12 14147.25  # L134: 11575.6362 -> 14147.25
Set-CellValue $ws 134 13 -107147715  # M134: -107148105 -> -107147715
Set-CellValue $ws 134 14 -19217.25  # N134: -16645.6362 -> -19217.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
Set-CellValue $ws 16 8 57045.777  # H16: 64249.875 -> 57045.777
Set-CellValue $ws 16 9 84437  # I16: 101240 -> 84437
Set-CellValue $ws 16 10 2263.3333  # J16: 2599.6667 -> 2263.3333
Set-CellValue $ws 16 11 84437  # K16: 101240 -> 84437
Set-CellValue $ws 16 12 2263.3333  # L16: 2599.6667 -> 2263.3333
Set-CellValue $ws 16 13 -84150  # M16: -100953 -> -84150
Set-CellValue $ws 16 14 -2837.3333  # N16: -3173.6667 -> -2837.3333
# Row 31
Set-CellValue $ws 31 8 2465.25  # H31: 2399.606 -> 2465.25
Set-CellValue $ws 31 9 1435.9474  # I31: 1361.0476 -> 1435.9474
Set-CellValue $ws 31 10 3969.6155  # J31: 4217.0835 -> 3969.6155
Set-CellValue $ws 31 11 1435.9474  # K31: 1361.0476 -> 1435.9474
Set-CellValue $ws 31 12 3969.6155  # L31: 4217.0835 -> 3969.6155
Set-CellValue $ws 31 13 -1140.9474  # M31: -1066.0476 -> -1140.9474
Set-CellValue $ws 31 14 -4559.6155  # N31: -4807.0835 -> -4559.6155
# Row 34
Set-CellValue $ws 34 8 2465.25  # H34: 2399.606 -> 2465.25
Set-CellValue $ws 34 9 1435.9474  # I34: 1361.0476 -> 1435.9474
Set-CellValue $ws 34 10 3969.6155  # J34: 4217.0835 -> 3969.6155
Set-CellValue $ws 34 11 1435.9474  # K34: 1361.0476 -> 1435.9474
Set-CellValue $ws 34 12 3969.6155  # L34: 4217.0835 -> 3969.6155
Set-CellValue $ws 34 13 -1233.9474  # M34: -1159.0476 -> -1233.9474
Set-CellValue $ws 34 14 -4373.6155  # N34: -4621.0835 -> -4373.6155
# Row 58
Set-CellValue $ws 58 8 1647.8889  # H58: 1477.475 -> 1647.8889
Set-CellValue $ws 58 9 1089.2593  # I58: 963.63635 -> 1089.2593
Set-CellValue $ws 58 10 3323.7778  # J58: 3899.8572 -> 3323.7778
Set-CellValue $ws 58 11 1089.2593  # K58: 963.63635 -> 1089.2593
Set-CellValue $ws 58 12 3323.7778  # L58: 3899.8572 -> 3323.7778
Set-CellValue $ws 58 13 -886.2592999999999  # M58: -760.63635 -> -886.2592999999999
Set-CellValue $ws 58 14 -3729.7778  # N58: -4305.8572 -> -3729.7778
# Row 99
Set-CellValue $ws 99 8 12517574  # H99: 20860128 -> 12517574
Set-CellValue $ws 99 9 15646222  # I99: 31289444 -> 15646222
Set-CellValue $ws 99 10 2980  # J99: 1499 -> 2980
Set-CellValue $ws 99 11 15646222  # K99: 31289444 -> 15646222
Set-CellValue $ws 99 12 2980  # L99: 1499 -> 2980
Set-CellValue $ws 99 13 -15644724  # M99: -31287946 -> -15644724
Set-CellValue $ws 99 14 -5976  # N99: -4495 -> -5976
# Row 113
Set-CellValue $ws 113 8 57045.777  # H113: 64249.875 -> 57045.777
Set-CellValue $ws 113 9 84437  # I113: 101240 -> 84437
Set-CellValue $ws 113 10 2263.3333  # J113: 2599.6667 -> 2263.3333
Set-CellValue $ws 113 11 84437  # K113: 101240 -> 84437
Set-CellValue $ws 113 12 2263.3333  # L113: 2599.6667 -> 2263.3333
Set-CellValue $ws 113 13 -82267  # M113: -99070 -> -82267
Set-CellValue $ws 113 14 -6603.3333  # N113: -6939.6667 -> -6603.3333
# Row 126
Set-CellValue $ws 126 8 12517574  # H126: 20860128 -> 12517574
Set-CellValue $ws 126 9 15646222  # I126: 31289444 -> 15646222
Set-CellValue $ws 126 10 2980  # J126: 1499 -> 2980
Set-CellValue $ws 126 11 46938666  # K126: 93868332 -> 46938666
Set-CellValue $ws 126 12 8940  # L126: 4497 -> 8940
Set-CellValue $ws 126 13 -46936196  # M126: -93865862 -> -46936196
Set-CellValue $ws 126 14 -13880  # N126: -9437 -> -13880
# Row 132
Set-CellValue $ws 132 8 3164.4412  # H132: 2164.6086 -> 3164.4412
Set-CellValue $ws 132 9 2690.2917  # I132: 1818.3 -> 2690.2917
Set-CellValue $ws 132 10 4302.4  # J132: 4473.3335 -> 4302.4
Set-CellValue $ws 132 11 8070.875100000001  # K132: 5454.9 -> 8070.875100000001
Set-CellValue $ws 132 12 12907.2  # L132: 13420.0005 -> 12907.2
Set-CellValue $ws 132 13 -5540.875100000001  # M132: -2924.9 -> -5540.875100000001
Set-CellValue $ws 132 14 -17967.2  # N132: -18480.0005 -> -17967.2
# Row 134
Set-CellValue $ws 134 8 3027.6177  # H134: 1796.0615 -> 3027.6177
Set-CellValue $ws 134 9 1944.6154  # I134: 1217.1818 -> 1944.6154
Set-CellValue $ws 134 10 6547.375  # J134: 4979.9 -> 6547.375
Set-CellValue $ws 134 11 5833.8462  # K134: 3651.5454 -> 5833.8462
Set-CellValue $ws 134 12 19642.125  # L134: 14939.7 -> 19642.125
Set-CellValue $ws 134 13 -3298.8462  # M134: -1116.5454 -> -3298.8462
Set-CellValue $ws 134 14 -24712.125  # N134: -20009.7 -> -24712.125
# Row 136
Set-CellValue $ws 136 8 1647.8889  # H136: 1477.475 -> 1647.8889
Set-CellValue $ws 136 9 1089.2593  # I136: 963.63635 -> 1089.2593
Set-CellValue $ws 136 10 3323.7778  # J136: 3899.8572 -> 3323.7778
Set-CellValue $ws 136 11 3267.7779  # K136: 2890.90905 -> 3267.7779
Set-CellValue $ws 136 12 9971.3334  # L136: 11699.5716 -> 9971.3334
Set-CellValue $ws 136 13 -717.7779  # M136: -340.9090500000002 -> -717.7779
Set-CellValue $ws 136 14 -15071.3334  # N136: -16799.5716 -> -15071.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
Set-CellValue $ws 131 8 7093995.5  # H131: 6668412.5 -> 7093995.5
Set-CellValue $ws 131 10 8335359  # J131: 7753888.5 -> 8335359
Set-CellValue $ws 131 12 25006077  # L131: 23261665.5 -> 25006077
Set-CellValue $ws 131 14 -25016157  # N131: -23271745.5 -> -25016157

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
Set-CellValue $ws 113 8 1992.6923  # H113: 1934.3928 -> 1992.6923
Set-CellValue $ws 113 9 1551.6923  # I113: 1466.9412 -> 1551.6923
Set-CellValue $ws 113 10 2433.6924  # J113: 2656.818 -> 2433.6924
Set-CellValue $ws 113 11 1551.6923  # K113: 1466.9412 -> 1551.6923
Set-CellValue $ws 113 12 2433.6924  # L113: 2656.818 -> 2433.6924
Set-CellValue $ws 113 13 618.3077000000001  # M113: 703.0588 -> 618.3077000000001
Set-CellValue $ws 113 14 -6773.6924  # N113: -6996.818 -> -6773.6924
# Row 122
Set-CellValue $ws 122 8 742670.4  # H122: 696215.0600000001 -> 742670.4
Set-CellValue $ws 122 9 1236204.5  # I122: 1011621.2 -> 1236204.5
Set-CellValue $ws 122 10 2369.3333  # J122: 2321.6 -> 2369.3333
Set-CellValue $ws 122 11 3708613.5  # K122: 3034863.6 -> 3708613.5
Set-CellValue $ws 122 12 7107.999899999999  # L122: 6964.799999999999 -> 7107.999899999999
Set-CellValue $ws 122 13 -3706163.5  # M122: -3032413.6 -> -3706163.5
Set-CellValue $ws 122 14 -12007.9999  # N122: -11864.8 -> -12007.9999
# Row 132
Set-CellValue $ws 132 8 3740.122  # H132: 2990.7593 -> 3740.122
Set-CellValue $ws 132 9 3972.9583  # I132: 2905.1143 -> 3972.9583
Set-CellValue $ws 132 10 3411.4119  # J132: 3148.5264 -> 3411.4119
Set-CellValue $ws 132 11 11918.8749  # K132: 8715.3429 -> 11918.8749
Set-CellValue $ws 132 12 10234.2357  # L132: 9445.5792 -> 10234.2357
Set-CellValue $ws 132 13 -9388.874899999999  # M132: -6185.3429 -> -9388.874899999999
Set-CellValue $ws 132 14 -15294.2357  # N132: -14505.5792 -> -15294.2357
# Row 138
Set-CellValue $ws 138 8 78000  # H138: 0 -> 78000
Set-CellValue $ws 138 10 78000  # J138: 0 -> 78000
Set-CellValue $ws 138 12 78000  # L138: 0 -> 78000
Set-CellValue $ws 138 14 -88280  # N138: None -> -88280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
Set-CellValue $ws 7 8 3220.04  # H7: 3494 -> 3220.04
Set-CellValue $ws 7 9 2589.5557  # I7: 2968 -> 2589.5557
Set-CellValue $ws 7 10 3574.6875  # J7: 3615.3845 -> 3574.6875
Set-CellValue $ws 7 11 2589.5557  # K7: 2968 -> 2589.5557
Set-CellValue $ws 7 12 3574.6875  # L7: 3615.3845 -> 3574.6875
Set-CellValue $ws 7 13 -2477.5557  # M7: -2856 -> -2477.5557
Set-CellValue $ws 7 14 -3798.6875  # N7: -3839.3845 -> -3798.6875
# Row 22
Set-CellValue $ws 22 8 6667  # H22: 7626.4707 -> 6667
Set-CellValue $ws 22 9 1076  # I22: 1216.6666 -> 1076
Set-CellValue $ws 22 10 8530.666999999999  # J22: 9000 -> 8530.666999999999
Set-CellValue $ws 22 11 1076  # K22: 1216.6666 -> 1076
Set-CellValue $ws 22 12 8530.666999999999  # L22: 9000 -> 8530.666999999999
Set-CellValue $ws 22 13 -781  # M22: -921.6666 -> -781
Set-CellValue $ws 22 14 -9120.666999999999  # N22: -9590 -> -9120.666999999999
# Row 27
Set-CellValue $ws 27 8 6667  # H27: 7626.4707 -> 6667
Set-CellValue $ws 27 9 1076  # I27: 1216.6666 -> 1076
Set-CellValue $ws 27 10 8530.666999999999  # J27: 9000 -> 8530.666999999999
Set-CellValue $ws 27 11 1076  # K27: 1216.6666 -> 1076
Set-CellValue $ws 27 12 8530.666999999999  # L27: 9000 -> 8530.666999999999
Set-CellValue $ws 27 13 -969  # M27: -1109.6666 -> -969
Set-CellValue $ws 27 14 -8744.666999999999  # N27: -9214 -> -8744.666999999999
# Row 126
Set-CellValue $ws 126 8 3220.04  # H126: 3494 -> 3220.04
Set-CellValue $ws 126 9 2589.5557  # I126: 2968 -> 2589.5557
Set-CellValue $ws 126 10 3574.6875  # J126: 3615.3845 -> 3574.6875
Set-CellValue $ws 126 11 7768.6671  # K126: 8904 -> 7768.6671
Set-CellValue $ws 126 12 10724.0625  # L126: 10846.1535 -> 10724.0625
Set-CellValue $ws 126 13 -5298.6671  # M126: -6434 -> -5298.6671
Set-CellValue $ws 126 14 -15664.0625  # N126: -15786.1535 -> -15664.0625
# Row 132
Set-CellValue $ws 132 8 6409.6763  # H132: 4014.1304 -> 6409.6763
Set-CellValue $ws 132 9 6021.913  # I132: 2047 -> 6021.913
Set-CellValue $ws 132 10 7220.4546  # J132: 7702.5 -> 7220.4546
Set-CellValue $ws 132 11 18065.739  # K132: 6141 -> 18065.739
Set-CellValue $ws 132 12 21661.3638  # L132: 23107.5 -> 21661.3638
Set-CellValue $ws 132 13 -15535.739  # M132: -3611 -> -15535.739
Set-CellValue $ws 132 14 -26721.3638  # N132: -28167.5 -> -26721.3638
# Row 136
Set-CellValue $ws 136 8 4223.6304  # H136: 3566.4814 -> 4223.6304
Set-CellValue $ws 136 9 2243.889  # I136: 1886.289 -> 2243.889
Set-CellValue $ws 136 10 11350.7  # J136: 11967.444 -> 11350.7
Set-CellValue $ws 136 11 6731.667  # K136: 5658.867 -> 6731.667
Set-CellValue $ws 136 12 34052.10000000001  # L136: 35902.33199999999 -> 34052.10000000001
Set-CellValue $ws 136 13 -4181.667  # M136: -3108.867 -> -4181.667
Set-CellValue $ws 136 14 -39152.10000000001  # N136: -41002.33199999999 -> -39152.10000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
Set-CellValue $ws 122 8 126437.875  # H122: 44990.477 -> 126437.875
Set-CellValue $ws 122 9 201000.6  # I122: 72758.64 -> 201000.6
Set-CellValue $ws 122 10 2166.6667  # J122: 1795.5555 -> 2166.6667
Set-CellValue $ws 122 11 603001.8  # K122: 218275.92 -> 603001.8
Set-CellValue $ws 122 12 6500.000100000001  # L122: 5386.666499999999 -> 6500.000100000001
Set-CellValue $ws 122 13 -600551.8  # M122: -215825.92 -> -600551.8
Set-CellValue $ws 122 14 -11400.0001  # N122: -10286.6665 -> -11400.0001
# Row 126
Set-CellValue $ws 126 8 51169.65  # H126: 51361.4 -> 51169.65
Set-CellValue $ws 126 9 63587.062  # I126: 72610.57000000001 -> 63587.062
Set-CellValue $ws 126 10 1500  # J126: 1780 -> 1500
Set-CellValue $ws 126 11 190761.186  # K126: 217831.71 -> 190761.186
Set-CellValue $ws 126 12 4500  # L126: 5340 -> 4500
Set-CellValue $ws 126 13 -188291.186  # M126: -215361.71 -> -188291.186
Set-CellValue $ws 126 14 -9440  # N126: -10280 -> -9440
# Row 132
Set-CellValue $ws 132 8 7145401.5  # H132: 1040.84 -> 7145401.5
Set-CellValue $ws 132 9 9618073  # I132: 863 -> 9618073
Set-CellValue $ws 132 10 2129.3333  # J132: 1521.6666 -> 2129.3333
Set-CellValue $ws 132 11 28854219  # K132: 2589 -> 28854219
Set-CellValue $ws 132 12 6387.999899999999  # L132: 4564.9998 -> 6387.999899999999
Set-CellValue $ws 132 13 -28851689  # M132: -59 -> -28851689
Set-CellValue $ws 132 14 -11447.9999  # N132: -9624.9998 -> -11447.9999
# Row 136
Set-CellValue $ws 136 8 16115.909  # H136: 19321.637 -> 16115.909
Set-CellValue $ws 136 9 22789.934  # I136: 27703.648 -> 22789.934
Set-CellValue $ws 136 10 1814.4286  # J136: 2091.9443 -> 1814.4286
Set-CellValue $ws 136 11 68369.802  # K136: 83110.944 -> 68369.802
Set-CellValue $ws 136 12 5443.2858  # L136: 6275.8329 -> 5443.2858
Set-CellValue $ws 136 13 -65819.802  # M136: -80560.944 -> -65819.802
Set-CellValue $ws 136 14 -10543.2858  # N136: -11375.8329 -> -10543.2858
